$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update translation status column (C) for rows translated by "juminho".
# Every row where column C = "Não" and column D = "juminho" becomes "SIM",
# except row 468 which becomes "sIM".
for ($r = 1; $r -le 782; $r++) {
    $status = $ws.Cells.Item($r, 3).Value
    $translator = $ws.Cells.Item($r, 4).Value
    if ($status -eq "Não" -and $translator -eq "juminho") {
        if ($r -eq 468) {
            $ws.Cells.Item($r, 3).Value = "sIM"
        } else {
            $ws.Cells.Item($r, 3).Value = "SIM"
        }
    }
}

# Update the view state (scroll position / active selection) to match.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 670
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C692").Select()
